# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells whose new values look numeric to remain
# text (matching the existing text-based Price column formatting)
$textPriceCells = @("D4", "D5", "D6", "D9", "D10", "D11", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Update cell values
$ws.Range("D2").Value = "60.580.29"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.318.26"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "559.48"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").Value = "142.80"
$ws.Range("E6").Value = "  -3.50%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.318.99"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Value = "7.90"
$ws.Range("E10").Value = "  -1.72%  "
$ws.Range("D11").Value = "0.119"
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").Value = "0.410"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "3.881.17"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "27.03"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").Value = "3.306.20"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "0.0000166"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "60.502.93"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "6.22"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "14.47"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "8.71"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").Value = "374.96"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "74.55"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.541"
$ws.Range("E25").Value = "  -4.42%  "
$ws.Range("D26").Value = "3.441.52"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("E27").Value = "  -6.37%  "
$ws.Range("D28").Value = "0.173"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "7.24"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "7.64"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").Value = "2.05"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D34").Value = "22.67"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "1.25"
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("D36").Value = "5.19"
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "165.98"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.75"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "26.95"
$ws.Range("E40").Value = "  -12.92%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.343.76"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "0.0736"
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D43").Value = "42.10"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "0.754"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").Value = "2.372.07"
$ws.Range("E48").Value = "  -6.47%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("D51").Value = "21.43"
$ws.Range("E51").Value = "  -4.92%  "
